# "changed stock to holding"
# The last header cell (E1) currently reads "Stock"; rename it to "Holding".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Holding"

# The body of the sheet below the header row (A2:E10) also picks up an
# explicit solid white interior fill as part of this change.
$ws.Range("A2:E10").Interior.Color = 16777215
